# Update 'want to go' (想去人数) counts in F column across sheets
# as captured at a later scrape run (gh-pages output regeneration).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 478
$ws.Range("F3").Value = 165
$ws.Range("F4").Value = 8090
$ws.Range("F7").Value = 1285
$ws.Range("F10").Value = 487
$ws.Range("F11").Value = 171
$ws.Range("F12").Value = 18
$ws.Range("F13").Value = 462
$ws.Range("F15").Value = 83
$ws.Range("F17").Value = 6020
$ws.Range("F18").Value = 196
$ws.Range("F19").Value = 287
$ws.Range("F20").Value = 2128
$ws.Range("F21").Value = 69
$ws.Range("F22").Value = 110
$ws.Range("F23").Value = 243
$ws.Range("F24").Value = 432

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 478
$ws.Range("F3").Value = 165
$ws.Range("F4").Value = 8090
$ws.Range("F7").Value = 1285
$ws.Range("F11").Value = 487
$ws.Range("F12").Value = 171
$ws.Range("F13").Value = 18
$ws.Range("F14").Value = 462
$ws.Range("F16").Value = 83
$ws.Range("F18").Value = 4
$ws.Range("F19").Value = 6020
$ws.Range("F21").Value = 196
$ws.Range("F22").Value = 287
$ws.Range("F23").Value = 2128
$ws.Range("F24").Value = 69
$ws.Range("F25").Value = 110
$ws.Range("F26").Value = 243
$ws.Range("F27").Value = 432
